$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: #US47 -> description changed to "Crear documento de testing" ---
$ws.Range("A3").Value2 = "#US47 Crear documento de testing"
$ws.Range("B3").Value2 = 25
$ws.Range("C3").Formula = "=IF(B3<SUM(F3:BL3),SUM(F3:BL3),B3)"
$ws.Range("D3").Formula = '=IF(C3>B3,$C3-(SUM($F3:$BL3)),$B3-(SUM($F3:$BL3)))'
$ws.Range("E3").Value2 = 0
$ws.Range("F3").Value2 = 5
$ws.Range("G3").Value2 = 2
$ws.Range("H3").Value2 = 5
$ws.Range("I3").Value2 = 3
$ws.Range("L3").Value2 = 10
$ws.Range("M3").Value2 = 0
$ws.Range("N3").Value2 = 0
$ws.Range("O3").Value2 = 0
$ws.Range("P3").Value2 = 0
$ws.Range("S3").Value2 = 0
$ws.Range("T3").Value2 = 0
$ws.Range("U3").Value2 = 0
$ws.Range("V3").Value2 = 0
$ws.Range("W3").Value2 = 0

# --- Row 4: #US97 Agregar funcionalidad al boton imprimir ---
$ws.Range("B4").Value2 = 20
$ws.Range("E4").Value2 = 10
$ws.Range("F4").Value2 = 0
$ws.Range("G4").Value2 = 10
$ws.Range("H4").Value2 = 0
$ws.Range("I4").Value2 = 0
$ws.Range("L4").Value2 = 0
$ws.Range("M4").Value2 = 0
$ws.Range("N4").Value2 = 0
$ws.Range("O4").Value2 = 0
$ws.Range("P4").Value2 = 0
$ws.Range("S4").Value2 = 0
$ws.Range("T4").Value2 = 0
$ws.Range("U4").Value2 = 0
$ws.Range("V4").Value2 = 0
$ws.Range("W4").Value2 = 0

# --- Row 5: #US98 Sacar Boton visualizar de interfaz ---
$ws.Range("B5").Value2 = 20
$ws.Range("E5").Value2 = 0
$ws.Range("F5").Value2 = 0
$ws.Range("G5").Value2 = 0
$ws.Range("H5").Value2 = 0
$ws.Range("I5").Value2 = 0
$ws.Range("L5").Value2 = 0
$ws.Range("M5").Value2 = 0
$ws.Range("N5").Value2 = 0
$ws.Range("O5").Value2 = 0
$ws.Range("P5").Value2 = 0
$ws.Range("S5").Value2 = 0
$ws.Range("T5").Value2 = 5
$ws.Range("U5").Value2 = 10
$ws.Range("V5").Value2 = 3
$ws.Range("W5").Value2 = 2

# --- Row 6: #US99 Detenet Red Ad-Hoc ---
$ws.Range("B6").Value2 = 25
$ws.Range("E6").Value2 = 0
$ws.Range("F6").Value2 = 5
$ws.Range("G6").Value2 = 0
$ws.Range("H6").Value2 = 10
$ws.Range("I6").Value2 = 0
$ws.Range("L6").Value2 = 0
$ws.Range("M6").Value2 = 10
$ws.Range("N6").Value2 = 0
$ws.Range("O6").Value2 = 0
$ws.Range("P6").Value2 = 0
$ws.Range("S6").Value2 = 0
$ws.Range("T6").Value2 = 0
$ws.Range("U6").Value2 = 0
$ws.Range("V6").Value2 = 0
$ws.Range("W6").Value2 = 0

# --- Row 7: #US102 Mostrar mensaje de error cuando alumno pone mal ip ---
$ws.Range("E7").Value2 = 0
$ws.Range("F7").Value2 = 0
$ws.Range("G7").Value2 = 0
$ws.Range("H7").Value2 = 0
$ws.Range("I7").Value2 = 0
$ws.Range("L7").Value2 = 0
$ws.Range("M7").Value2 = 0
$ws.Range("N7").Value2 = 0
$ws.Range("O7").Value2 = 5
$ws.Range("P7").Value2 = 0
$ws.Range("S7").Value2 = 10
$ws.Range("T7").Value2 = 0
$ws.Range("U7").Value2 = 0
$ws.Range("V7").Value2 = 0
$ws.Range("W7").Value2 = 0

# --- Update the last-selected cell in the bottom-right frozen pane ---
$ws.Range("X5").Select()
